$d = $word.ActiveDocument

$replacements = @(
    @{old = "97×13=1261"; new = "63×65=4095"},
    @{old = "70×46=3220"; new = "85×25=2125"},
    @{old = "29×77=2233"; new = "98×28=2744"},
    @{old = "43×66=2838"; new = "86×44=3784"},
    @{old = "12×54=648";  new = "25×13=325"},
    @{old = "18×84=1512"; new = "77×95=7315"},
    @{old = "17×76=1292"; new = "35×84=2940"},
    @{old = "58×52=3016"; new = "49×65=3185"},
    @{old = "20×42=840";  new = "78×35=2730"},
    @{old = "95×95=9025"; new = "26×77=2002"},
    @{old = "54×79=4266"; new = "84×53=4452"},
    @{old = "60×37=2220"; new = "56×36=2016"},
    @{old = "99×12=1188"; new = "65×60=3900"},
    @{old = "56×27=1512"; new = "63×77=4851"},
    @{old = "55×39=2145"; new = "77×75=5775"},
    @{old = "81×14=1134"; new = "17×64=1088"},
    @{old = "62×99=6138"; new = "87×85=7395"},
    @{old = "53×70=3710"; new = "35×43=1505"},
    @{old = "40×87=3480"; new = "71×79=5609"},
    @{old = "39×70=2730"; new = "61×54=3294"},
    @{old = "64×93=5952"; new = "65×54=3510"},
    @{old = "15×86=1290"; new = "98×40=3920"},
    @{old = "99×17=1683"; new = "56×51=2856"},
    @{old = "94×62=5828"; new = "82×31=2542"},
    @{old = "64×23=1472"; new = "42×35=1470"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
